# Applies the GTA 6 Implementation Program SOW rewrite described in the commit diff.
$d = $word.ActiveDocument
$lb = [char]11   # vertical tab -> Word line break (<w:br/>), matches the w:br separators in the diff

# Paragraph 2: Title heading
$d.Paragraphs.Item(2).Range.Text = "GTA 6 Implementation Program"

# Paragraph 3: Opening SOW paragraph (parties / dates)
$d.Paragraphs.Item(3).Range.Text = "This Statement of Work (“SOW”) is entered into as of [INSERT SOW EFFECTIVE DATE] by and between The Service Provider is a leading technology and software development organization with a proven track record in large-scale, multi-platform project implementations. The company maintains robust quality assurance, security, and compliance processes and operates globally across multiple time zones, applying industry best practices in agile and waterfall project management. (“[INSERT SERVICE PROVIDER NAME]”) and The Client is a recognized enterprise with strategic interest in advancing next-generation gaming environments and interactive digital experiences. The Client is committed to innovation, quality, and timely release schedules, and maintains an experienced in-house team to collaborate with the Service Provider. (“[INSERT CLIENT NAME]”) under the provisions of that certain Master Services Agreement, dated as of [INSERT AGREEMENT DATE], by and between [INSERT SERVICE PROVIDER NAME] and [INSERT CLIENT NAME] (the “Agreement”)."

# Paragraph 5: Services Description body
$d.Paragraphs.Item(5).Range.Text = "The Service Provider shall deliver a full-cycle implementation of GTA 6, including requirements analysis, architecture, design, development, integration, testing, deployment, and post-deployment hypercare. Methodologies will include:  " + $lb + "- Agile and Waterfall hybrid project management to optimize flexibility and delivery predictability  " + $lb + "- Detailed requirements gathering with iterative validation and sign-off  " + $lb + "- High-fidelity architectural design covering both technical and business domains  " + $lb + "- Secure, scalable, and maintainable coding standards and rigorous code reviews  " + $lb + "- Comprehensive test planning, including unit, integration, system, performance, and security testing  " + $lb + "- Continuous integration and deployment pipelines to ensure rapid feedback and quality control  " + $lb + "- Structured knowledge transfer and operational handover  " + $lb + "- Ongoing stakeholder communication, risk management, and issue escalation mechanisms  " + $lb + "- All services will be performed in accordance with relevant ISO, SOC, and other applicable certification standards, ensuring data integrity and regulatory compliance throughout the engagement."

# Paragraph 7: Deliverables body
$d.Paragraphs.Item(7).Range.Text = "- Business Requirements Specification Document: Comprehensive documentation of all technical, functional, and non-functional requirements for GTA 6.  " + $lb + "- Systems Architecture Design: Detailed architecture blueprints covering infrastructure, networking, data models, security layers, and scalability provisions.  " + $lb + "- UI/UX Design Prototypes: High-fidelity interface and user journey mockups, adhering to accessibility and usability standards.  " + $lb + "- Core Application Source Code: Fully functional, documented, and version-controlled codebase for the GTA 6 platform.  " + $lb + "- Integration Modules: Source code, API documentation, and test cases for all required external and internal system integrations.  " + $lb + "- Automated Test Suites: Scripts and documentation for unit, integration, system, and regression testing.  " + $lb + "- Deployment Artifacts: All binaries, configuration files, and environment setup scripts necessary for production deployment.  " + $lb + "- User and Technical Documentation: Complete manuals, operational guides, and training materials for end-users and support staff.  " + $lb + "- Post-Go-Live Support Plan: Structured hypercare support processes, escalation paths, and service desk arrangements.  " + $lb + "- Project Status and Financial Reporting: Regular progress, risk, and financial status reports, ensuring transparency and client oversight."

# Paragraph 9: Milestones body
$d.Paragraphs.Item(9).Range.Text = "- Milestone 1: Completion of requirements gathering and approval of Business Requirements Specification Document ([INSERT DATE])  " + $lb + "- Milestone 2: Delivery and sign-off of Systems Architecture Design ([INSERT DATE])  " + $lb + "- Milestone 3: Finalization and client approval of UI/UX Design Prototypes ([INSERT DATE])  " + $lb + "- Milestone 4: Completion of core application development and first full-system integration ([INSERT DATE])  " + $lb + "- Milestone 5: Completion of system, performance, and security testing ([INSERT DATE])  " + $lb + "- Milestone 6: Production deployment and operational handover ([INSERT DATE])  " + $lb + "- Milestone 7: Completion of post-go-live hypercare and project closure ([INSERT DATE])"

# Paragraph 11: Acceptance body
$d.Paragraphs.Item(11).Range.Text = "Acceptance of all deliverables and milestones will be based on the following criteria:  " + $lb + "- Formal sign-off by the Client’s designated representative(s) upon review of submitted deliverables  " + $lb + "- Compliance with the technical specifications, business requirements, and quality standards defined in this SOW  " + $lb + "- Successful completion of user acceptance testing (UAT) for all functional components  " + $lb + "- Resolution of all critical defects identified during review or testing phases  " + $lb + "- Written notification by the Client of acceptance or a detailed list of deficiencies within ten (10) business days of deliverable submission. If deficiencies are reported, the Service Provider will remedy and resubmit corrected deliverables within five (5) business days. Deliverables will be deemed accepted if no feedback is issued within the review period."

# Paragraph 13: Personnel and Locations body
$d.Paragraphs.Item(13).Range.Text = "- Primary work will be performed at the Service Provider’s primary development center and approved secure remote locations.  " + $lb + "- Key personnel will include Project Manager, Solution Architect, Lead Developers, QA/Test Manager, Integration Engineer(s), UI/UX Designer(s), and Technical Writer(s).  " + $lb + "- All personnel assigned to the project will have relevant qualifications, certifications, and experience as required for their respective roles.  " + $lb + "- Occasional on-site collaboration at the Client’s location(s) may be scheduled by mutual agreement, subject to health, safety, and security protocols."

# Paragraph 14: Representatives heading
$d.Paragraphs.Item(14).Range.Text = "Representatives"

# Paragraph 15: Representatives body
$d.Paragraphs.Item(15).Range.Text = "- Project Manager: [INSERT NAME]  " + $lb + "- Solution Architect: [INSERT NAME]  " + $lb + "- Lead Developer: [INSERT NAME]  " + $lb + "- QA/Test Manager: [INSERT NAME]  " + $lb + "- Key contact details for escalation and decision-making will be provided upon project commencement."

# Paragraph 17: Client Representatives body
$d.Paragraphs.Item(17).Range.Text = "- Program Sponsor: [INSERT NAME]  " + $lb + "- IT Director: [INSERT NAME]  " + $lb + "- Solution Owner: [INSERT NAME]  " + $lb + "- Technical Point of Contact: [INSERT NAME]  " + $lb + "- Additional stakeholders may be designated as needed."

# Paragraph 18: Contractor Resources heading
$d.Paragraphs.Item(18).Range.Text = "Contractor Resources"

# Paragraph 19: Contractor Resources body
$d.Paragraphs.Item(19).Range.Text = "- The Service Provider will allocate sufficient qualified resources to ensure successful delivery, including:  " + $lb + "  - 1 Project Manager (full-time)  " + $lb + "  - 1 Solution Architect (full-time)  " + $lb + "  - Development team (front-end, back-end, integration specialists, database engineers) sized appropriately for each project phase  " + $lb + "  - QA/testing specialists for automated and manual testing  " + $lb + "  - UI/UX designers and technical writers  " + $lb + "- All personnel will be subject to background and compliance checks where required by the Client.  " + $lb + "- Staffing adjustments may be made by mutual agreement to address project needs."

# Paragraph 20: Terms & Conditions heading
$d.Paragraphs.Item(20).Range.Text = "Terms & Conditions"

# Paragraph 21: Terms & Conditions body
$d.Paragraphs.Item(21).Range.Text = "This Statement of Work shall commence on the Start Date and will remain in effect until completion of all deliverables and milestones or until the End Date, unless terminated earlier in accordance with the Termination clause.  " + $lb + "- The total project duration is estimated at twelve (12) years from the Start Date, as detailed in the Timeline section.  " + $lb + "- Adjustments to key deadlines and milestones may be made only through the Change Process outlined herein.  " + $lb + "- The SOW Effective Date and Agreement Date are as specified in the header sections and will be confirmed upon contract execution."

# Paragraph 23: Fees body
$d.Paragraphs.Item(23).Range.Text = "This engagement is structured as a fixed-price contract. The total fixed fee for all services, deliverables, and milestones described in this SOW shall be agreed upon prior to commencement and documented in the signed contract.  " + $lb + "- This fee covers all labor, project management, technical services, design, development, testing, documentation, and deployment activities as specified herein.  " + $lb + "- Any work outside the defined scope will require prior written authorization and may be subject to additional charges at mutually agreed rates."

# Paragraph 25: Expenses body
$d.Paragraphs.Item(25).Range.Text = "The Service Provider shall be reimbursed for pre-approved, reasonable out-of-pocket expenses incurred in direct connection with the performance of the services under this SOW.  " + $lb + "- Such expenses may include travel, lodging, subsistence, and necessary materials, and will be invoiced at actual cost with supporting receipts.  " + $lb + "- All such expenses must be pre-approved in writing by the Client and are subject to the Client’s expense policies."

# Paragraph 27: Taxes body
$d.Paragraphs.Item(27).Range.Text = "All fees and expenses are exclusive of any applicable federal, state, local, or international taxes, levies, duties, or similar governmental assessments.  " + $lb + "- The Client shall be responsible for all taxes arising from or in connection with the services provided under this SOW, except for taxes based on the Service Provider’s income, employment, or property."

# Paragraph 29: Conversion body
$d.Paragraphs.Item(29).Range.Text = "If the Client elects to engage any member of the Service Provider’s project team as a direct employee or contractor within twelve (12) months of their last involvement in this project, the Client shall pay a conversion fee equal to twenty-five percent (25%) of the individual’s annualized compensation. This fee shall be due within thirty (30) days of the individual’s hiring date."

# Paragraph 31: Limitation of Liability body
$d.Paragraphs.Item(31).Range.Text = "To the fullest extent permitted by law, neither party shall be liable to the other for any indirect, incidental, consequential, special, or punitive damages, including loss of profits, data, or business opportunities, arising out of or related to this SOW, regardless of the cause of action.  " + $lb + "- The aggregate liability of either party for direct damages under this SOW shall not exceed the total fees paid or payable by the Client to the Service Provider under this SOW during the twelve (12) month period preceding the event giving rise to liability.  " + $lb + "- Any claims under this SOW must be brought within three (3) months after discovery or sixty (60) days following SOW termination, whichever is earlier."

# Paragraph 33: Service Level Agreement body
$d.Paragraphs.Item(33).Range.Text = "- The Service Provider will deliver all services and deliverables in a professional and workmanlike manner, consistent with industry best practices.  " + $lb + "- All project milestones will be met as scheduled, with any anticipated deviations communicated in writing within five (5) business days of identification.  " + $lb + "- System uptime, performance, and security standards will comply with the specifications set forth in the design and architecture documentation.  " + $lb + "- Post-production support will provide response to critical incidents within four (4) business hours and resolution within agreed service windows.  " + $lb + "- Regular monitoring, reporting, and continuous improvement processes will be implemented."

# Paragraph 35: Assumptions body
$d.Paragraphs.Item(35).Range.Text = "- The Client will provide timely access to all necessary systems, environments, documentation, and key personnel.  " + $lb + "- All third-party dependencies, software licenses, and hardware required for the project will be made available as scheduled.  " + $lb + "- The project scope is as defined in this SOW; any deviations or additions will be handled through the Change Process.  " + $lb + "- No legacy system data migration or support for out-of-scope platforms (e.g., GTA 5 PC) is included.  " + $lb + "- The Service Provider is not responsible for delays or cost overruns resulting from factors beyond its reasonable control, including force majeure events."

# Paragraph 37: Change Process body
$d.Paragraphs.Item(37).Range.Text = "Any changes to the project scope, schedule, deliverables, or fees must be requested in writing and processed through a formal Change Request procedure.  " + $lb + "- Each Change Request shall document the nature of the change, its impact on the project, and any adjustments to cost or timeline.  " + $lb + "- No change shall be implemented without written approval from authorized representatives of both parties.  " + $lb + "- Approved changes shall be appended to this SOW and become binding upon execution."

# Paragraph 40: Client signature line
$d.Paragraphs.Item(40).Range.Text = "[INSERT CLIENT NAME]  Signature: ________________"

# Paragraph 41: Service Provider signature line
$d.Paragraphs.Item(41).Range.Text = "[INSERT SERVICE PROVIDER NAME] Signature: ________________"
